# "Klik Advance System Settings"
#
# The last paragraph of the document (a ListParagraph bullet item) holds
# only a floating (anchored) screenshot picture and a trailing tab
# character; the paragraph also hosts the document's hidden "_GoBack"
# bookmark (the last-edit marker Word maintains automatically), which
# currently wraps the picture run.
#
# The edit:
#   1. Appends a tab run + a new text run ("Advance System Settings")
#      after the existing picture run (the pre-existing trailing tab
#      run ends up immediately followed by the new text).
#   2. Relocates the "_GoBack" bookmark so it collapses to a single
#      point right after the new text (i.e. at the very end of the
#      paragraph, before the paragraph mark) instead of wrapping the
#      picture run.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$paraRange = $lastPara.Range

# 1) Add the new tab + text content at the end of the paragraph (after
#    the picture run and the pre-existing tab run that already live
#    there).
$paraRange.InsertAfter("Advance System Settings")

# 2) Move the "_GoBack" bookmark to collapse right before the paragraph
#    mark, i.e. immediately after the text just inserted.
#
#    Bookmark placement right at the very last position(s) of the
#    document's Content range is unreliable, so a scratch character is
#    appended past the real end-of-document first to keep the desired
#    bookmark position safely away from that boundary; it is removed
#    again afterwards via the same Range reference (no text search
#    involved, so nothing elsewhere in the document can be touched).
$targetPos = $d.Content.End - 1

$scratchRange = $d.Range($d.Content.End, $d.Content.End)
$scratchRange.InsertAfter("x")

$bookmarkRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$scratchRange.Delete()
